$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item(1)

# The row labeled "grandes regiões e unidades da federação" (row 6, a
# header-only row with no data) was removed; everything below it shifts
# up one row so the numeric data lines up correctly with the region
# labels (e.g. "norte" now carries the values that used to belong to
# "rondônia", etc. — i.e. the whole block slides up by one).
$ws.Rows.Item(6).Delete()
